# Apply cell updates from the crypto price/volume refresh.
# Values that look like plain numbers (e.g. "620.97") are written with a
# leading apostrophe (Excel quote-prefix) so they stay text, matching how
# the Price column is stored in the source workbook (inline text, not numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '90.485.93'
$ws.Range('E2').Value = '  +2.95%  '
$ws.Range('D3').Value = '3.194.04'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '''620.97'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('D7').Value = '''0.402'
$ws.Range('E7').Value = '  +4.71%  '
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '3.187.41'
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('E12').Value = '  -6.12%  '
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('D14').Value = '90.171.33'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('D15').Value = '3.783.21'
$ws.Range('E15').Value = '  -1.72%  '
$ws.Range('E16').Value = '  -2.65%  '
$ws.Range('D17').Value = '''5.27'
$ws.Range('E17').Value = '  -3.80%  '
$ws.Range('D18').Value = '3.181.16'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('D19').Value = '''3.28'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').Value = '''0.0000213'
$ws.Range('E20').Value = '  +49.87%  '
$ws.Range('D21').Value = '''13.49'
$ws.Range('E21').Value = '  -3.74%  '
$ws.Range('D22').Value = '''440.05'
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('E23').Value = '  -3.50%  '
$ws.Range('E24').Value = '  -4.52%  '
$ws.Range('D25').Value = '''5.17'
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('E26').Value = '  -5.98%  '
$ws.Range('D27').Value = '3.353.93'
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('D28').Value = '''75.57'
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '''0.170'
$ws.Range('E30').Value = '  -3.61%  '
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  +28.92%  '
$ws.Range('D33').Value = '''8.49'
$ws.Range('E33').Value = '  -3.88%  '
$ws.Range('D34').Value = '''538.99'
$ws.Range('E34').Value = '  -5.30%  '
$ws.Range('D35').Value = '''7.03'
$ws.Range('E35').Value = '  -3.06%  '
$ws.Range('E36').Value = '  -4.29%  '
$ws.Range('E37').Value = '  -8.11%  '
$ws.Range('D38').Value = '''22.13'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').Value = '''22.38'
$ws.Range('E39').Value = '  +2.62%  '
$ws.Range('E40').Value = '  -7.94%  '
$ws.Range('D41').Value = '''0.998'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('E42').Value = '  -3.72%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  -6.57%  '
$ws.Range('D45').Value = '''150.60'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''173.26'
$ws.Range('E46').Value = '  -3.57%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '''43.73'
$ws.Range('E47').Value = '  -3.28%  '
$ws.Range('E48').Value = '  -8.03%  '
$ws.Range('E49').Value = '  -7.55%  '
$ws.Range('D50').Value = '''4.08'
$ws.Range('E50').Value = '  -3.94%  '
$ws.Range('D51').Value = '''0.613'
$ws.Range('E51').Value = '  -3.18%  '
